$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.045.15'
$ws.Range('E2').Value = '  -0.50%  '
$ws.Range('D3').Value = '2.340.84'
$ws.Range('E3').Value = '  +1.49%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '307.12'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.15%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '101.71'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.81%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.510'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -4.31%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  -2.93%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '34.94'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -1.97%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '52.58'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +1.80%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0798'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -1.89%  '
$ws.Range('E13').Value = '  +0.98%  '
$ws.Range('E14').Value = '  -2.04%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '15.83'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +5.75%  '
$ws.Range('D16').Value = '2.343.78'
$ws.Range('E16').Value = '  +3.70%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.828'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +2.72%  '
$ws.Range('D18').Value = '42.988.56'
$ws.Range('E18').Value = '  -0.46%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.21'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.64%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '11.76'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -4.29%  '
$ws.Range('D21').Value = '0.0₃0912'
$ws.Range('E21').Value = '  -1.87%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '68.11'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '237.01'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -1.85%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.02'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.45%  '
$ws.Range('E25').Value = '  -1.81%  '
$ws.Range('E26').Value = '  -0.09%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '25.70'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +3.53%  '
$ws.Range('B28').Value = 'LEO'
$ws.Range('C28').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '3.95'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.57%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.32'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +4.56%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '35.82'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -2.04%  '
$ws.Range('B31').Value = 'Cosmos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '9.32'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -3.17%  '
$ws.Range('B32').Value = 'Monero'
$ws.Range('C32').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '161.90'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -4.39%  '
$ws.Range('B33').Value = 'FirstDigitalUSD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.999'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -0.15%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.12'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -2.79%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.67'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +8.36%  '
$ws.Range('B36').Value = 'Celestia'
$ws.Range('C36').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '17.54'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -1.04%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.0727'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -1.86%  '
$ws.Range('B38').Value = 'WEMIXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.44'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -3.97%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.86'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -1.23%  '
$ws.Range('B40').Value = 'LidoDAOToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.92'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -4.55%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.102'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -3.01%  '
$ws.Range('B42').Value = 'Stellar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.113'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -2.58%  '
$ws.Range('B43').Value = 'ApeXProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.61'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +10.59%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '2.020.83'
$ws.Range('E44').Value = '  +2.22%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0285'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -2.59%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '18.92'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.96%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '10.18'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +2.02%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.94'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -1.34%  '
$ws.Range('B49').Value = 'MultiversX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '56.47'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +2.21%  '
$ws.Range('B50').Value = 'HuobiToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.88'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -2.46%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '2.564.92'
$ws.Range('E51').Value = '  +1.12%  '
